$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (sheet 1) ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F3").Value = 171
$ws1.Range("F4").Value = 8677
$ws1.Range("F5").Value = 106
$ws1.Range("F9").Value = 117
$ws1.Range("F10").Value = 493
$ws1.Range("F11").Value = 180
$ws1.Range("F13").Value = 471
$ws1.Range("F16").Value = 34
$ws1.Range("F17").Value = 6159
$ws1.Range("F18").Value = 205
$ws1.Range("F19").Value = 307
$ws1.Range("F20").Value = 2287
$ws1.Range("F21").Value = 101
$ws1.Range("F22").Value = 171
$ws1.Range("F23").Value = 245
$ws1.Range("F24").Value = 455

# ---- Sheet "演出" (sheet 2) ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 8

# ---- Sheet "全部类型" (sheet 4) ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F3").Value = 171
$ws4.Range("F4").Value = 8677
$ws4.Range("F5").Value = 106
$ws4.Range("F9").Value = 8
$ws4.Range("F11").Value = 117
$ws4.Range("F12").Value = 493
$ws4.Range("F13").Value = 180
$ws4.Range("F15").Value = 471
$ws4.Range("F18").Value = 34
$ws4.Range("F20").Value = 6159
$ws4.Range("F22").Value = 205
$ws4.Range("F23").Value = 307
$ws4.Range("F24").Value = 2287
$ws4.Range("F25").Value = 101
$ws4.Range("F26").Value = 171
$ws4.Range("F27").Value = 245
$ws4.Range("F28").Value = 455
